$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows (dates as Excel serials, nuovi pos., somma mobile 7gg., somma mobile per 100mila ab.)
$data = @(
    @(344, 44418, 0, 2, 18.42468908337172),
    @(345, 44419, 0, 2, 18.42468908337172),
    @(346, 44420, 0, 2, 18.42468908337172),
    @(347, 44421, 3, 4, 36.84937816674343),
    @(348, 44422, 0, 4, 36.84937816674343),
    @(349, 44423, 3, 6, 55.27406725011516),
    @(350, 44424, 0, 6, 55.27406725011516),
    @(351, 44425, 2, 8, 73.69875633348687),
    @(352, 44426, 0, 8, 73.69875633348687),
    @(353, 44427, 0, 8, 73.69875633348687),
    @(354, 44428, 1, 6, 55.27406725011516),
    @(355, 44429, 0, 6, 55.27406725011516),
    @(356, 44430, 5, 8, 73.69875633348687),
    @(357, 44431, 0, 8, 73.69875633348687)
)

foreach ($row in $data) {
    $r = $row[0]
    $dateSerial = $row[1]
    $newPos = $row[2]
    $sumMobile = $row[3]
    $sumMobile100k = $row[4]

    # Copy the date-cell formatting (border/font/alignment/number format) from the
    # last existing date cell, then overwrite its value.
    $ws.Range("A343").Copy() | Out-Null
    $ws.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($r, 1).Value = $dateSerial

    $ws.Cells.Item($r, 2).Value = $newPos
    $ws.Cells.Item($r, 3).Value = $sumMobile
    $ws.Cells.Item($r, 4).Value = $sumMobile100k
}

$excel.CutCopyMode = 0
